# Generate Report for Handoff
#
# A new handoff round was generated: the "Status" cells that previously
# showed the handback state now show that the file is ready to be handed
# off again, the "latest xliff generate / handoff" timestamps advance to
# the new generation time, and the over-wide "...Datetime" columns are
# narrowed down to a more reasonable width on all three sheets.

$wb = $excel.ActiveWorkbook

# Target display width ~17.22 chars. ColumnWidth is quantized to whole
# pixels by Excel, so 16.3 is the input that lands on the closest
# representable column width.
$newWidth = 16.3

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-20 17:06:56"
$ws.Columns("E").ColumnWidth = $newWidth
$ws.Columns("F").ColumnWidth = $newWidth

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-20 17:06:52"
$ws.Columns("C").ColumnWidth = $newWidth

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-20 17:06:56"
$ws.Columns("C").ColumnWidth = $newWidth
